$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows/cols where the new value looks like a plain number and must be forced to text
# (matching the source data, which stores all D/E values as text strings).
$textRows = @(5, 6, 7, 10, 11, 12, 13, 16, 17, 21, 23, 25, 26, 27, 28, 29, 31, 32, 33, 34, 35, 37, 38, 39, 41, 42, 43, 44, 45, 46, 47, 48)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "43.025.77"
$ws.Cells.Item(2, 5).Value = "  +4.21%  "

$ws.Cells.Item(3, 4).Value = "2.261.29"
$ws.Cells.Item(3, 5).Value = "  +3.12%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).Value = "253.63"
$ws.Cells.Item(5, 5).Value = "  -0.61%  "

$ws.Cells.Item(6, 4).Value = "0.631"
$ws.Cells.Item(6, 5).Value = "  +1.96%  "

$ws.Cells.Item(7, 4).Value = "71.84"
$ws.Cells.Item(7, 5).Value = "  +4.59%  "

$ws.Cells.Item(8, 5).Value = "  +18.48%  "

$ws.Cells.Item(9, 5).Value = "  -0.01%  "

$ws.Cells.Item(10, 4).Value = "39.99"
$ws.Cells.Item(10, 5).Value = "  +6.43%  "

$ws.Cells.Item(11, 4).Value = "0.0983"
$ws.Cells.Item(11, 5).Value = "  +4.21%  "

$ws.Cells.Item(12, 4).Value = "59.56"
$ws.Cells.Item(12, 5).Value = "  +1.15%  "

$ws.Cells.Item(13, 4).Value = "7.60"
$ws.Cells.Item(13, 5).Value = "  +6.95%  "

$ws.Cells.Item(14, 5).Value = "  -0.12%  "

$ws.Cells.Item(15, 4).Value = "2.602.92"
$ws.Cells.Item(15, 5).Value = "  +3.50%  "

$ws.Cells.Item(16, 4).Value = "0.889"
$ws.Cells.Item(16, 5).Value = "  +1.07%  "

$ws.Cells.Item(17, 4).Value = "14.85"
$ws.Cells.Item(17, 5).Value = "  +2.37%  "

$ws.Cells.Item(18, 4).Value = "2.262.89"
$ws.Cells.Item(18, 5).Value = "  +4.37%  "

$ws.Cells.Item(19, 4).Value = "42.948.64"
$ws.Cells.Item(19, 5).Value = "  +4.13%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0988"
$ws.Cells.Item(20, 5).Value = "  +2.86%  "

$ws.Cells.Item(21, 4).Value = "6.31"
$ws.Cells.Item(21, 5).Value = "  +1.88%  "

$ws.Cells.Item(22, 5).Value = "  +1.51%  "

$ws.Cells.Item(23, 4).Value = "236.17"
$ws.Cells.Item(23, 5).Value = "  +0.97%  "

$ws.Cells.Item(24, 5).Value = "  +4.86%  "

$ws.Cells.Item(25, 4).Value = "3.96"
$ws.Cells.Item(25, 5).Value = "  +2.03%  "

$ws.Cells.Item(26, 4).Value = "11.60"
$ws.Cells.Item(26, 5).Value = "  -2.27%  "

$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  -0.04%  "

$ws.Cells.Item(28, 4).Value = "2.46"
$ws.Cells.Item(28, 5).Value = "  -1.80%  "

$ws.Cells.Item(29, 4).Value = "3.67"
$ws.Cells.Item(29, 5).Value = "  -1.35%  "

$ws.Cells.Item(30, 5).Value = "  +8.47%  "

$ws.Cells.Item(31, 4).Value = "168.14"
$ws.Cells.Item(31, 5).Value = "  -0.67%  "

$ws.Cells.Item(32, 4).Value = "21.35"
$ws.Cells.Item(32, 5).Value = "  +3.12%  "

$ws.Cells.Item(33, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(33, 4).Value = "6.37"
$ws.Cells.Item(33, 5).Value = "  +15.70%  "

$ws.Cells.Item(34, 2).Value = "Kaspa"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(34, 4).Value = "0.129"
$ws.Cells.Item(34, 5).Value = "  +10.02%  "

$ws.Cells.Item(35, 4).Value = "0.0775"
$ws.Cells.Item(35, 5).Value = "  +2.53%  "

$ws.Cells.Item(36, 5).Value = "  +2.04%  "

$ws.Cells.Item(37, 4).Value = "29.42"
$ws.Cells.Item(37, 5).Value = "  +13.45%  "

$ws.Cells.Item(38, 4).Value = "4.74"
$ws.Cells.Item(38, 5).Value = "  +2.61%  "

$ws.Cells.Item(39, 4).Value = "4.17"
$ws.Cells.Item(39, 5).Value = "  +0.43%  "

$ws.Cells.Item(40, 5).Value = "  +8.27%  "

$ws.Cells.Item(41, 4).Value = "2.31"
$ws.Cells.Item(41, 5).Value = "  +4.13%  "

$ws.Cells.Item(42, 4).Value = "5.87"
$ws.Cells.Item(42, 5).Value = "  +3.53%  "

$ws.Cells.Item(43, 4).Value = "12.37"
$ws.Cells.Item(43, 5).Value = "  +1.46%  "

$ws.Cells.Item(44, 4).Value = "64.52"
$ws.Cells.Item(44, 5).Value = "  +1.37%  "

$ws.Cells.Item(45, 4).Value = "5.03"
$ws.Cells.Item(45, 5).Value = "  +1.43%  "

$ws.Cells.Item(46, 4).Value = "0.202"
$ws.Cells.Item(46, 5).Value = "  +1.76%  "

$ws.Cells.Item(47, 4).Value = "8.96"
$ws.Cells.Item(47, 5).Value = "  +3.00%  "

$ws.Cells.Item(48, 4).Value = "0.103"
$ws.Cells.Item(48, 5).Value = "  +1.04%  "

$ws.Cells.Item(49, 5).Value = "  -4.72%  "

$ws.Cells.Item(50, 5).Value = "  -0.10%  "

$ws.Cells.Item(51, 5).Value = "  +2.04%  "
